# Update cryptos list: prices (D) and 1h volume/change (E) columns,
# plus a Cronos/HuobiToken row swap (B/C/D/E) in rows 47-48.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '39.443.16'
$ws.Cells.Item(2, 5).Value = '  +1.70%  '

$ws.Cells.Item(3, 4).Value = '2.160.66'
$ws.Cells.Item(3, 5).Value = '  +3.24%  '

$ws.Cells.Item(4, 5).Value = '  +0.05%  '

$ws.Cells.Item(5, 4).Value = '''228.41'
$ws.Cells.Item(5, 5).Value = '  -0.29%  '

$ws.Cells.Item(6, 5).Value = '  +0.99%  '

$ws.Cells.Item(7, 4).Value = '''64.30'
$ws.Cells.Item(7, 5).Value = '  +5.05%  '

$ws.Cells.Item(8, 5).Value = '  +0.02%  '

$ws.Cells.Item(9, 4).Value = '''0.398'
$ws.Cells.Item(9, 5).Value = '  +3.03%  '

$ws.Cells.Item(10, 5).Value = '  +2.07%  '

$ws.Cells.Item(11, 5).Value = '  +0.40%  '

$ws.Cells.Item(12, 5).Value = '  +4.21%  '

$ws.Cells.Item(13, 4).Value = '2.482.62'
$ws.Cells.Item(13, 5).Value = '  +3.30%  '

$ws.Cells.Item(14, 4).Value = '''22.34'
$ws.Cells.Item(14, 5).Value = '  +1.31%  '

$ws.Cells.Item(15, 5).Value = '  +1.17%  '

$ws.Cells.Item(16, 4).Value = '''5.55'
$ws.Cells.Item(16, 5).Value = '  +0.98%  '

$ws.Cells.Item(17, 4).Value = '2.160.75'
$ws.Cells.Item(17, 5).Value = '  +1.98%  '

$ws.Cells.Item(18, 4).Value = '39.352.86'
$ws.Cells.Item(18, 5).Value = '  +1.66%  '

$ws.Cells.Item(19, 4).Value = '''71.87'
$ws.Cells.Item(19, 5).Value = '  +0.06%  '

$ws.Cells.Item(20, 4).Value = '''6.12'
$ws.Cells.Item(20, 5).Value = '  +0.75%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0857'
$ws.Cells.Item(21, 5).Value = '  +1.96%  '

$ws.Cells.Item(23, 5).Value = '  +0.12%  '

$ws.Cells.Item(24, 4).Value = '''2.51'
$ws.Cells.Item(24, 5).Value = '  +6.02%  '

$ws.Cells.Item(25, 4).Value = '''2.37'
$ws.Cells.Item(25, 5).Value = '  +1.25%  '

$ws.Cells.Item(26, 4).Value = '''172.30'
$ws.Cells.Item(26, 5).Value = '  +0.51%  '

$ws.Cells.Item(27, 4).Value = '''9.53'
$ws.Cells.Item(27, 5).Value = '  -0.11%  '

$ws.Cells.Item(28, 5).Value = '  +0.87%  '

$ws.Cells.Item(29, 4).Value = '''19.91'
$ws.Cells.Item(29, 5).Value = '  +2.99%  '

$ws.Cells.Item(30, 5).Value = '  -1.35%  '

$ws.Cells.Item(31, 4).Value = '''2.64'
$ws.Cells.Item(31, 5).Value = '  +5.82%  '

$ws.Cells.Item(32, 5).Value = '  +1.34%  '

$ws.Cells.Item(33, 5).Value = '  +2.42%  '

$ws.Cells.Item(34, 4).Value = '''4.76'
$ws.Cells.Item(34, 5).Value = '  +0.15%  '

$ws.Cells.Item(35, 4).Value = '''7.07'
$ws.Cells.Item(35, 5).Value = '  +9.23%  '

$ws.Cells.Item(36, 4).Value = '''0.0621'
$ws.Cells.Item(36, 5).Value = '  +1.30%  '

$ws.Cells.Item(37, 5).Value = '  +0.23%  '

$ws.Cells.Item(38, 5).Value = '  +0.34%  '

$ws.Cells.Item(39, 5).Value = '  -0.07%  '

$ws.Cells.Item(40, 5).Value = '  +1.33%  '

$ws.Cells.Item(41, 4).Value = '''103.76'
$ws.Cells.Item(41, 5).Value = '  +2.89%  '

$ws.Cells.Item(42, 4).Value = '''17.82'
$ws.Cells.Item(42, 5).Value = '  -0.84%  '

$ws.Cells.Item(43, 4).Value = '1.540.03'
$ws.Cells.Item(43, 5).Value = '  +0.31%  '

$ws.Cells.Item(44, 5).Value = '  +4.13%  '

$ws.Cells.Item(45, 4).Value = '''7.96'
$ws.Cells.Item(45, 5).Value = '  +3.53%  '

$ws.Cells.Item(46, 5).Value = '  +4.30%  '

$ws.Cells.Item(47, 2).Value = 'HuobiToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(47, 4).Value = '''2.82'
$ws.Cells.Item(47, 5).Value = '  +0.56%  '

$ws.Cells.Item(48, 2).Value = 'Cronos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(48, 4).Value = '''0.0925'
$ws.Cells.Item(48, 5).Value = '  +1.38%  '

$ws.Cells.Item(49, 5).Value = '  +5.95%  '

$ws.Cells.Item(50, 4).Value = '2.365.25'
$ws.Cells.Item(50, 5).Value = '  +3.33%  '
